# Apply the crypto price/volume refresh described by the commit diff.
# Numeric-looking "Price" values are forced back to literal text (NumberFormat
# "@" + Style "Normal") so Excel does not silently coerce strings such as
# "1.000" or "0.07446" into floating point numbers (which would lose the
# original formatting/precision), while leaving the cell style unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.009.20'
$ws.Range('E2').Value = '  -1.37%  '
$ws.Range('D3').Value = '1.781.74'
$ws.Range('E3').Value = '  -1.45%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.17'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9999'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.10%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5392'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.17%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3771'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.33%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07446'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.97%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.65'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.37%  '
$ws.Range('E11').Value = '  -2.82%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9998'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.08%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.47'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.24%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.078'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.97%  '
$ws.Range('E15').Value = '  -1.78%  '
$ws.Range('D16').Value = '1.775.83'
$ws.Range('E16').Value = '  -1.57%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '88.30'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -4.36%  '
$ws.Range('E18').Value = '  -1.59%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06439'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9998'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.07%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.21'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.876'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.11%  '
$ws.Range('D23').Value = '28.027.99'
$ws.Range('E23').Value = '  -1.31%  '
$ws.Range('E24').Value = '  -2.09%  '
$ws.Range('E25').Value = '  -2.03%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '155.96'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.56%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.22'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.54%  '
$ws.Range('D28').Value = '1.979.19'
$ws.Range('E28').Value = '  -1.67%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.277'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.16%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '120.03'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.105'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.09%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1055'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.41%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.642'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.85%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.515'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.2253'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.91%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06440'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.00%  '
$ws.Range('E37').Value = '  -1.69%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.003'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.46%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.422'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.71%  '
$ws.Range('B40').Value = 'WEMIXTOKEN'
$ws.Range('C40').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.445'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.67%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6139'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.45%  '
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.07'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.177'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.58%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9993'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.20'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.39%  '
$ws.Range('E46').Value = '  -0.46%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5741'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.09%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '126.28'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.88%  '
$ws.Range('E49').Value = '  +3.23%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.922'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.44%  '
$ws.Range('E51').Value = '  -1.69%  '
